$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "login" to "LoginTest"
$ws.Name = "LoginTest"

# Remove the "address" column (column C) entirely
$ws.Columns.Item(3).Delete()

# Remove the last data row (row 5) entirely
$ws.Rows.Item(5).Delete()

# Update the selected cell/range shown in the worksheet view
$ws.Range("G9").Select()
